# Apply updated cryptocurrency market data (prices / 1h volume change)
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.132.85"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.319.38"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.57"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.07"
$ws.Range("E10").Value = "  +4.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.70"
$ws.Range("E13").Value = "  -2.56%  "

$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.680.60"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.317.46"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("E17").Value = "  -3.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.026.89"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +5.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.22"
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.72"
$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("E24").Value = "  -2.75%  "

$ws.Range("E25").Value = "  -0.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.19"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.09"
$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.69"
$ws.Range("E31").Value = "  -1.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.96"
$ws.Range("E32").Value = "  +5.11%  "

$ws.Range("E33").Value = "  +2.70%  "

$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.41"
$ws.Range("E34").Value = "  +7.17%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("E36").Value = "  -0.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0695"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("E38").Value = "  +1.24%  "

$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("E40").Value = "  -1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.998.44"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.16"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.52"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.85"
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "76.95"
$ws.Range("E47").Value = "  +9.10%  "

$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.86"
$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.545.91"
$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  +3.35%  "
